$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new transaction row is inserted above the current row 2, pushing the
# existing rows 2-4 down to rows 3-5. Write the cells directly (bottom row
# first) so each write lands on its final destination without relying on a
# native row-insert (which would also drag unwanted formatting along).

# old row 4 -> row 5
$ws.Range("E5").Value = "Deposit"
$ws.Range("N5").Value = "Crypto"
$ws.Range("P5").Value = "ETH"
$ws.Range("T5").Value = 341.28

# old row 3 -> row 4
$ws.Range("E4").Value = "Deposit"
$ws.Range("N4").Value = "Crypto"
$ws.Range("P4").Value = "ETH"
$ws.Range("T4").Value = 596.35320000000002

# old row 2 -> row 3
$ws.Range("E3").Value = "Withdrawal"
$ws.Range("N3").Value = "Wiretransfer"
$ws.Range("P3").Value = "Anywires"
$ws.Range("T3").Value = 1622.46

# new row 2
$ws.Range("E2").Value = "Withdrawal"
$ws.Range("N2").Value = "Crypto"
$ws.Range("P2").Value = "ETH"
$ws.Range("T2").Value = 500.02510000000001

# Matches the selection recorded in the saved workbook.
$ws.Range("I12").Select()
